$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data. Values that look like numbers
# (e.g. "245.00", "36.688.11") must stay as literal text, so the
# touched range is switched to text format before assignment and
# restored to the default "Normal" style afterwards so no stray
# number-format/style attribute is introduced.

# Row 2
$cells = $ws.Range("D2:E2")
$cells.NumberFormat = "@"
$ws.Range("D2").Value = "36.688.11"
$ws.Range("E2").Value = "  +0.75%  "
$cells.Style = "Normal"

# Row 3
$cells = $ws.Range("D3")
$cells.NumberFormat = "@"
$ws.Range("D3").Value = "1.965.34"
$cells.Style = "Normal"

# Row 4
$cells = $ws.Range("E4")
$cells.NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$cells.Style = "Normal"

# Row 5
$cells = $ws.Range("D5:E5")
$cells.NumberFormat = "@"
$ws.Range("D5").Value = "245.00"
$ws.Range("E5").Value = "  +1.22%  "
$cells.Style = "Normal"

# Row 6
$cells = $ws.Range("D6:E6")
$cells.NumberFormat = "@"
$ws.Range("D6").Value = "0.619"
$ws.Range("E6").Value = "  +1.92%  "
$cells.Style = "Normal"

# Row 7
$cells = $ws.Range("D7:E7")
$cells.NumberFormat = "@"
$ws.Range("D7").Value = "61.53"
$ws.Range("E7").Value = "  +8.81%  "
$cells.Style = "Normal"

# Row 8
$cells = $ws.Range("E8")
$cells.NumberFormat = "@"
$ws.Range("E8").Value = "  +0.04%  "
$cells.Style = "Normal"

# Row 9
$cells = $ws.Range("E9")
$cells.NumberFormat = "@"
$ws.Range("E9").Value = "  +5.68%  "
$cells.Style = "Normal"

# Row 10
$cells = $ws.Range("D10:E10")
$cells.NumberFormat = "@"
$ws.Range("D10").Value = "0.0800"
$ws.Range("E10").Value = "  -5.31%  "
$cells.Style = "Normal"

# Row 11
$cells = $ws.Range("E11")
$cells.NumberFormat = "@"
$ws.Range("E11").Value = "  +0.75%  "
$cells.Style = "Normal"

# Row 12
$cells = $ws.Range("D12:E12")
$cells.NumberFormat = "@"
$ws.Range("D12").Value = "14.33"
$ws.Range("E12").Value = "  +7.41%  "
$cells.Style = "Normal"

# Row 13
$cells = $ws.Range("D13:E13")
$cells.NumberFormat = "@"
$ws.Range("D13").Value = "22.08"
$ws.Range("E13").Value = "  +4.55%  "
$cells.Style = "Normal"

# Row 14
$cells = $ws.Range("D14:E14")
$cells.NumberFormat = "@"
$ws.Range("D14").Value = "0.841"
$ws.Range("E14").Value = "  +4.83%  "
$cells.Style = "Normal"

# Row 15
$cells = $ws.Range("D15:E15")
$cells.NumberFormat = "@"
$ws.Range("D15").Value = "2.235.62"
$ws.Range("E15").Value = "  +0.86%  "
$cells.Style = "Normal"

# Row 16
$cells = $ws.Range("D16:E16")
$cells.NumberFormat = "@"
$ws.Range("D16").Value = "5.32"
$ws.Range("E16").Value = "  +4.61%  "
$cells.Style = "Normal"

# Row 17
$cells = $ws.Range("D17:E17")
$cells.NumberFormat = "@"
$ws.Range("D17").Value = "1.958.70"
$ws.Range("E17").Value = "  +1.07%  "
$cells.Style = "Normal"

# Row 18
$cells = $ws.Range("D18:E18")
$cells.NumberFormat = "@"
$ws.Range("D18").Value = "36.698.69"
$ws.Range("E18").Value = "  +0.96%  "
$cells.Style = "Normal"

# Row 19
$cells = $ws.Range("D19:E19")
$cells.NumberFormat = "@"
$ws.Range("D19").Value = "70.05"
$ws.Range("E19").Value = "  +1.77%  "
$cells.Style = "Normal"

# Row 20
$cells = $ws.Range("D20:E20")
$cells.NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0857"
$ws.Range("E20").Value = "  -0.16%  "
$cells.Style = "Normal"

# Row 21
$cells = $ws.Range("B21:E21")
$cells.NumberFormat = "@"
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "5.11"
$ws.Range("E21").Value = "  +3.20%  "
$cells.Style = "Normal"

# Row 22
$cells = $ws.Range("B22:E22")
$cells.NumberFormat = "@"
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "230.90"
$ws.Range("E22").Value = "  +2.02%  "
$cells.Style = "Normal"

# Row 23
$cells = $ws.Range("E23")
$cells.NumberFormat = "@"
$ws.Range("E23").Value = "  +0.06%  "
$cells.Style = "Normal"

# Row 24
$cells = $ws.Range("D24:E24")
$cells.NumberFormat = "@"
$ws.Range("D24").Value = "2.48"
$ws.Range("E24").Value = "  +7.42%  "
$cells.Style = "Normal"

# Row 25
$cells = $ws.Range("E25")
$cells.NumberFormat = "@"
$ws.Range("E25").Value = "  +3.71%  "
$cells.Style = "Normal"

# Row 26
$cells = $ws.Range("E26")
$cells.NumberFormat = "@"
$ws.Range("E26").Value = "  +7.11%  "
$cells.Style = "Normal"

# Row 27
$cells = $ws.Range("D27:E27")
$cells.NumberFormat = "@"
$ws.Range("D27").Value = "9.25"
$ws.Range("E27").Value = "  +2.28%  "
$cells.Style = "Normal"

# Row 28
$cells = $ws.Range("D28:E28")
$cells.NumberFormat = "@"
$ws.Range("D28").Value = "160.85"
$ws.Range("E28").Value = "  +0.07%  "
$cells.Style = "Normal"

# Row 29
$cells = $ws.Range("D29:E29")
$cells.NumberFormat = "@"
$ws.Range("D29").Value = "19.50"
$ws.Range("E29").Value = "  +1.77%  "
$cells.Style = "Normal"

# Row 30
$cells = $ws.Range("D30:E30")
$cells.NumberFormat = "@"
$ws.Range("D30").Value = "1.25"
$ws.Range("E30").Value = "  +12.31%  "
$cells.Style = "Normal"

# Row 31
$cells = $ws.Range("E31")
$cells.NumberFormat = "@"
$ws.Range("E31").Value = "  +2.18%  "
$cells.Style = "Normal"

# Row 32
$cells = $ws.Range("D32:E32")
$cells.NumberFormat = "@"
$ws.Range("D32").Value = "4.81"
$ws.Range("E32").Value = "  +6.52%  "
$cells.Style = "Normal"

# Row 33
$cells = $ws.Range("D33:E33")
$cells.NumberFormat = "@"
$ws.Range("D33").Value = "0.0621"
$ws.Range("E33").Value = "  +0.70%  "
$cells.Style = "Normal"

# Row 34
$cells = $ws.Range("D34:E34")
$cells.NumberFormat = "@"
$ws.Range("D34").Value = "4.48"
$ws.Range("E34").Value = "  +8.51%  "
$cells.Style = "Normal"

# Row 35
$cells = $ws.Range("D35:E35")
$cells.NumberFormat = "@"
$ws.Range("D35").Value = "3.53"
$ws.Range("E35").Value = "  +16.91%  "
$cells.Style = "Normal"

# Row 36
$cells = $ws.Range("E36")
$cells.NumberFormat = "@"
$ws.Range("E36").Value = "  +6.11%  "
$cells.Style = "Normal"

# Row 37
$cells = $ws.Range("E37")
$cells.NumberFormat = "@"
$ws.Range("E37").Value = "  -0.09%  "
$cells.Style = "Normal"

# Row 38
$cells = $ws.Range("E38")
$cells.NumberFormat = "@"
$ws.Range("E38").Value = "  -0.99%  "
$cells.Style = "Normal"

# Row 39
$cells = $ws.Range("D39")
$cells.NumberFormat = "@"
$ws.Range("D39").Value = "5.59"
$cells.Style = "Normal"

# Row 40
$cells = $ws.Range("D40:E40")
$cells.NumberFormat = "@"
$ws.Range("D40").Value = "0.0987"
$ws.Range("E40").Value = "  +0.31%  "
$cells.Style = "Normal"

# Row 41
$cells = $ws.Range("E41")
$cells.NumberFormat = "@"
$ws.Range("E41").Value = "  +1.09%  "
$cells.Style = "Normal"

# Row 42
$cells = $ws.Range("E42")
$cells.NumberFormat = "@"
$ws.Range("E42").Value = "  +3.26%  "
$cells.Style = "Normal"

# Row 43
$cells = $ws.Range("E43")
$cells.NumberFormat = "@"
$ws.Range("E43").Value = "  +1.99%  "
$cells.Style = "Normal"

# Row 44
$cells = $ws.Range("D44:E44")
$cells.NumberFormat = "@"
$ws.Range("D44").Value = "16.23"
$ws.Range("E44").Value = "  +4.87%  "
$cells.Style = "Normal"

# Row 45
$cells = $ws.Range("D45:E45")
$cells.NumberFormat = "@"
$ws.Range("D45").Value = "1.369.84"
$ws.Range("E45").Value = "  +2.80%  "
$cells.Style = "Normal"

# Row 46
$cells = $ws.Range("B46:E46")
$cells.NumberFormat = "@"
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "1.04"
$ws.Range("E46").Value = "  +2.75%  "
$cells.Style = "Normal"

# Row 47
$cells = $ws.Range("B47:E47")
$cells.NumberFormat = "@"
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "89.11"
$ws.Range("E47").Value = "  +4.48%  "
$cells.Style = "Normal"

# Row 48
$cells = $ws.Range("D48:E48")
$cells.NumberFormat = "@"
$ws.Range("D48").Value = "7.19"
$ws.Range("E48").Value = "  +1.68%  "
$cells.Style = "Normal"

# Row 49
$cells = $ws.Range("E49")
$cells.NumberFormat = "@"
$ws.Range("E49").Value = "  +0.81%  "
$cells.Style = "Normal"

# Row 50
$cells = $ws.Range("D50:E50")
$cells.NumberFormat = "@"
$ws.Range("D50").Value = "44.64"
$ws.Range("E50").Value = "  +3.20%  "
$cells.Style = "Normal"

# Row 51
$cells = $ws.Range("D51:E51")
$cells.NumberFormat = "@"
$ws.Range("D51").Value = "2.129.86"
$ws.Range("E51").Value = "  +1.00%  "
$cells.Style = "Normal"

